$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.714.09"
$ws.Range("E2").Value = "  +6.77%  "

$ws.Range("D3").Value = "'1.811.89"
$ws.Range("E3").Value = "  +4.93%  "

$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'250.72"
$ws.Range("E5").Value = "  +3.51%  "

$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "'0.4977"
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("D8").Value = "'0.2782"
$ws.Range("E8").Value = "  +7.36%  "

$ws.Range("D9").Value = "'0.06382"
$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("D10").Value = "'1.807.38"
$ws.Range("E10").Value = "  +4.66%  "

$ws.Range("D11").Value = "'16.72"

$ws.Range("D12").Value = "'0.07112"
$ws.Range("E12").Value = "  +3.00%  "

$ws.Range("D13").Value = "'0.6467"
$ws.Range("E13").Value = "  +6.31%  "

$ws.Range("D14").Value = "'4.699"
$ws.Range("E14").Value = "  +4.73%  "

$ws.Range("D15").Value = "'81.75"
$ws.Range("E15").Value = "  +5.81%  "

$ws.Range("D16").Value = "'28.687.24"
$ws.Range("E16").Value = "  +7.67%  "

$ws.Range("D17").Value = "'0.9990"

$ws.Range("D18").Value = "'0.000007380"
$ws.Range("E18").Value = "  +2.74%  "

$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").Value = "'12.27"
$ws.Range("E20").Value = "  +7.24%  "

$ws.Range("D21").Value = "'2.037.94"
$ws.Range("E21").Value = "  +4.33%  "

$ws.Range("E22").Value = "  +4.35%  "

$ws.Range("D23").Value = "'8.881"
$ws.Range("E23").Value = "  +3.67%  "

$ws.Range("D24").Value = "'5.321"
$ws.Range("E24").Value = "  +4.36%  "

$ws.Range("D25").Value = "'143.13"
$ws.Range("E25").Value = "  +3.35%  "

$ws.Range("D26").Value = "'16.03"
$ws.Range("E26").Value = "  +4.57%  "

$ws.Range("D27").Value = "'1.876"
$ws.Range("E27").Value = "  +5.27%  "

$ws.Range("D28").Value = "'112.77"
$ws.Range("E28").Value = "  +6.02%  "

$ws.Range("D29").Value = "'1.403"
$ws.Range("E29").Value = "  +1.61%  "

$ws.Range("D30").Value = "'4.178"
$ws.Range("E30").Value = "  +5.88%  "

$ws.Range("D31").Value = "'0.08358"
$ws.Range("E31").Value = "  +4.52%  "

$ws.Range("D32").Value = "'3.842"
$ws.Range("E32").Value = "  +4.20%  "

$ws.Range("D33").Value = "'0.04963"
$ws.Range("E33").Value = "  +9.71%  "

$ws.Range("D34").Value = "'1.087"
$ws.Range("E34").Value = "  +7.82%  "

$ws.Range("D35").Value = "'0.6757"
$ws.Range("E35").Value = "  +8.17%  "

$ws.Range("D36").Value = "'2.665"
$ws.Range("E36").Value = "  +2.66%  "

$ws.Range("D37").Value = "'2.708"
$ws.Range("E37").Value = "  +10.42%  "

$ws.Range("D38").Value = "'0.9593"
$ws.Range("E38").Value = "  +2.30%  "

$ws.Range("D39").Value = "'2.143"
$ws.Range("E39").Value = "  +4.30%  "

$ws.Range("D40").Value = "'0.01590"
$ws.Range("E40").Value = "  +5.91%  "

$ws.Range("D41").Value = "'5.980"
$ws.Range("E41").Value = "  +5.63%  "

$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").Value = "'101.05"
$ws.Range("E43").Value = "  +1.53%  "

$ws.Range("D44").Value = "'0.4107"
$ws.Range("E44").Value = "  +6.66%  "

$ws.Range("D45").Value = "'7.185"
$ws.Range("E45").Value = "  +4.53%  "

$ws.Range("E46").Value = "  +5.37%  "

$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("D48").Value = "'8.205"
$ws.Range("E48").Value = "  +3.70%  "

$ws.Range("D49").Value = "'31.41"
$ws.Range("E49").Value = "  +4.24%  "

$ws.Range("D50").Value = "'0.3627"
$ws.Range("E50").Value = "  +7.39%  "

$ws.Range("D51").Value = "'1.301"
$ws.Range("E51").Value = "  +5.46%  "
